$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.504.81"
$ws.Range("E2").Value = "  +0.78%  "

$ws.Range("D3").Value = "3.184.68"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.97%  "

$ws.Range("D8").Value = "3.185.32"
$ws.Range("E8").Value = "  -0.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.81%  "

$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.510"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("E13").Value = "  -2.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.86%  "

$ws.Range("D15").Value = "3.711.27"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").Value = "66.547.30"
$ws.Range("E16").Value = "  +1.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").Value = "3.187.14"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "512.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.82%  "

$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "511.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0894"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0421"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "

$ws.Range("E40").Value = "  +5.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.303"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.03%  "

$ws.Range("D43").Value = "0.0₃0679"
$ws.Range("E43").Value = "  +7.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.74%  "

$ws.Range("D46").Value = "2.855.90"
$ws.Range("E46").Value = "  -5.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.28%  "

$ws.Range("E48").Value = "  +3.92%  "

$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.07%  "
